# Completed SavingsAccount test plan: fill in the "Preconditions" (E),
# "Method Inputs" (F) and "Expected Result" (G) columns for the six test
# cases covering __init__ and get_service_charges, plus adjust the row
# heights that grew to fit the new wrapped text and restore the cursor /
# selection to where the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Test case 1 (row 7): __init__ succeeds, attributes set correctly ---
$ws.Range("E7").Value = 'None'
$ws.Range("F7").Value = 'account_number = 9483914, client_number = 345,                              balance = 1559.49,                         date_created = date(2024, 1, 1),  minimum_balance = 50.0'
$ws.Range("G7").Value = 'The savings_account instance is created successfully with the attributes correctly set.'

# --- Test case 2 (row 8): __init__ with invalid minimum_balance type ---
$ws.Range("E8").Value = 'None'
$ws.Range("F8").Value = 'account_number = 9483914, client_number = 345,                              balance = 1559.49,                         date_created = date(2024, 1, 1),  minimum_balance = "minimum"'
$ws.Range("G8").Value = 'minimum_balance = 50.0'

# --- Test case 3 (row 9): get_service_charges, balance greater than minimum ---
$ws.Range("E9").Value = 'account_number = 9483914, client_number = 345,                              balance = 1559.49,                         date_created = date(2024, 1, 1),  minimum_balance = 50.0'
$ws.Range("F9").Value = 'None'
$ws.Range("G9").Value = 'service_charge = BASE_SERVICE_CHARGE'

# --- Test case 4 (row 10): get_service_charges, balance equal to minimum ---
$ws.Range("E10").Value = 'account_number = 9483914, client_number = 345,                              balance = 50,                         date_created = date(2024, 1, 1),  minimum_balance = 50.0'
$ws.Range("F10").Value = 'None'
$ws.Range("G10").Value = 'service_charge = BASE_SERVICE_CHARGE'

# --- Test case 5 (row 11): get_service_charges, balance less than minimum ---
$ws.Range("E11").Value = 'account_number = 9483914, client_number = 345,                              balance = 10,                         date_created = date(2024, 1, 1),  minimum_balance = 50.0'
$ws.Range("F11").Value = 'None'
$ws.Range("G11").Value = 'service_charge = BASE_SERVICE_CHARGE *SERVICE_CHARGE_PREMIUM = 1.00'

# --- Test case 6 (row 12): __str__ ---
$ws.Range("E12").Value = 'account_number = 9483914, client_number = 345,                              balance = 1559.49,                         date_created = date(2024, 1, 1),  minimum_balance = 50.0'
$ws.Range("F12").Value = 'None'
$ws.Range("G12").Value = '"Account Number:9483914 Balance:$1,559.49\nMinimum Balance: $50.0 Account Type: Savings"'

# --- Row heights grew once the wrapped test data was typed in ---
$ws.Rows.Item(7).RowHeight = 75
$ws.Rows.Item(8).RowHeight = 72.75
$ws.Rows.Item(9).RowHeight = 70.9
$ws.Rows.Item(10).RowHeight = 83.25
$ws.Rows.Item(11).RowHeight = 103.25
$ws.Rows.Item(12).RowHeight = 78.4

# --- Restore the view: scrolled near the top, cursor left on F14 ---
$ws.Range("F14").Select()
